# Insert a new weekly record row at row 222, shifting existing rows 222-324 down to 223-325.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 222; this shifts rows 222..324 down to 223..325
$ws.Rows.Item(222).Insert()

# Populate the new row 222 with the new record's data
$ws.Cells.Item(222, 1).Value = 7
$ws.Cells.Item(222, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(222, 3).Value = "Ñuble"
$ws.Cells.Item(222, 4).Value = 44846
$ws.Cells.Item(222, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(222, 5).Value = 16
$ws.Cells.Item(222, 6).Value = 100112023
$ws.Cells.Item(222, 7).Value = "Brócoli"
$ws.Cells.Item(222, 8).Value = "Sin especificar"
$ws.Cells.Item(222, 9).Value = "Segunda"
$ws.Cells.Item(222, 10).Value = 200
$ws.Cells.Item(222, 11).Value = 1000
$ws.Cells.Item(222, 12).Value = 1000
$ws.Cells.Item(222, 13).Value = 1000
$ws.Cells.Item(222, 14).Value = "$/unidad"
$ws.Cells.Item(222, 15).Value = "Región del Maule"
$ws.Cells.Item(222, 16).Value = 1000
$ws.Cells.Item(222, 17).Value = 1
$ws.Cells.Item(222, 18).Value = "Hortaliza"
